$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1974.8026
$ws.Range("I15").Value = 1974.8026
$ws.Range("K15").Value = 5924.4078
$ws.Range("M15").Value = -5755.4078

$ws.Range("H33").Value = 1032.0714
$ws.Range("I33").Value = 889.95
$ws.Range("K33").Value = 889.95
$ws.Range("M33").Value = -660.95

$ws.Range("H64").Value = 3240
$ws.Range("I64").Value = 3300
$ws.Range("K64").Value = 3300
$ws.Range("M64").Value = -3052

$ws.Range("H67").Value = 3240
$ws.Range("I67").Value = 3300
$ws.Range("K67").Value = 3300
$ws.Range("M67").Value = -2442

$ws.Range("H76").Value = 3280.2
$ws.Range("I76").Value = 3183.5833
$ws.Range("J76").Value = 3666.6667
$ws.Range("K76").Value = 3183.5833
$ws.Range("L76").Value = 3666.6667
$ws.Range("M76").Value = -2868.5833
$ws.Range("N76").Value = -4296.6667

$ws.Range("H79").Value = 3280.2
$ws.Range("I79").Value = 3183.5833
$ws.Range("J79").Value = 3666.6667
$ws.Range("K79").Value = 3183.5833
$ws.Range("L79").Value = 3666.6667
$ws.Range("M79").Value = -2091.5833
$ws.Range("N79").Value = -5850.6667

$ws.Range("H113").Value = 2908.9
$ws.Range("I113").Value = 2877.8
$ws.Range("J113").Value = 2940
$ws.Range("K113").Value = 2877.8
$ws.Range("L113").Value = 2940
$ws.Range("M113").Value = 376.1999999999998
$ws.Range("N113").Value = -9448

$ws.Range("H132").Value = 4775.4644
$ws.Range("I132").Value = 4582
$ws.Range("K132").Value = 13746
$ws.Range("M132").Value = -11216

$ws.Range("H138").Value = 183249.69
$ws.Range("J138").Value = 280692.78
$ws.Range("L138").Value = 842078.3400000001
$ws.Range("N138").Value = -852358.3400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3308.3333
$ws.Range("J61").Value = 3985.7144
$ws.Range("L61").Value = 3985.7144
$ws.Range("N61").Value = -4409.7144

$ws.Range("H136").Value = 3308.3333
$ws.Range("J136").Value = 3985.7144
$ws.Range("L136").Value = 11957.1432
$ws.Range("N136").Value = -17057.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2425.372
$ws.Range("I134").Value = 2060.5483
$ws.Range("J134").Value = 3367.8333
$ws.Range("K134").Value = 6181.644899999999
$ws.Range("L134").Value = 10103.4999
$ws.Range("M134").Value = -3646.644899999999
$ws.Range("N134").Value = -15173.4999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 861.375
$ws.Range("I16").Value = 844.2
$ws.Range("K16").Value = 844.2
$ws.Range("M16").Value = -557.2

$ws.Range("H113").Value = 861.375
$ws.Range("I113").Value = 844.2
$ws.Range("K113").Value = 844.2
$ws.Range("M113").Value = 1325.8

$ws.Range("H132").Value = 11112884
$ws.Range("I132").Value = 1008.2
$ws.Range("K132").Value = 3024.6
$ws.Range("M132").Value = -494.6000000000004

$ws.Range("H134").Value = 1951.2727
$ws.Range("I134").Value = 1961.1111
$ws.Range("J134").Value = 1907
$ws.Range("K134").Value = 5883.3333
$ws.Range("L134").Value = 5721
$ws.Range("M134").Value = -3348.3333
$ws.Range("N134").Value = -10791

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1120.25
$ws.Range("J5").Value = 1367.25
$ws.Range("L5").Value = 4101.75
$ws.Range("N5").Value = -4325.75

$ws.Range("H39").Value = 1520.6842
$ws.Range("J39").Value = 1833
$ws.Range("L39").Value = 5499
$ws.Range("N39").Value = -6087

$ws.Range("H60").Value = 2797.9167
$ws.Range("I60").Value = 383.33334
$ws.Range("J60").Value = 3142.8572
$ws.Range("K60").Value = 1150.00002
$ws.Range("L60").Value = 9428.571599999999
$ws.Range("M60").Value = -899.0000199999999
$ws.Range("N60").Value = -9930.571599999999

$ws.Range("H69").Value = 37038376
$ws.Range("J69").Value = 55557056
$ws.Range("L69").Value = 166671168
$ws.Range("N69").Value = -166672790

$ws.Range("H72").Value = 37038376
$ws.Range("J72").Value = 55557056
$ws.Range("L72").Value = 500013504
$ws.Range("N72").Value = -500021616

$ws.Range("H110").Value = 11467.952
$ws.Range("I110").Value = 3006.75
$ws.Range("K110").Value = 9020.25
$ws.Range("M110").Value = -4930.25

$ws.Range("H128").Value = 200000
$ws.Range("I128").Value = 200000
$ws.Range("K128").Value = 600000
$ws.Range("M128").Value = -595020

$ws.Range("H134").Value = 2669.125
$ws.Range("I134").Value = 2618.111
$ws.Range("J134").Value = 2822.1667
$ws.Range("K134").Value = 7854.333
$ws.Range("L134").Value = 8466.500100000001
$ws.Range("M134").Value = -2784.333
$ws.Range("N134").Value = -18606.5001

$ws.Range("H135").Value = 1120.25
$ws.Range("J135").Value = 1367.25
$ws.Range("L135").Value = 12305.25
$ws.Range("N135").Value = -17375.25

$ws.Range("H140").Value = 3081.6667
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 3081.6667
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 9245.000100000001
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -19605.0001

$ws.Range("H141").Value = 6460.1
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 6460.1
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 19380.3
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -29740.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5837.143
$ws.Range("I70").Value = 5812.8
$ws.Range("J70").Value = 5862.5
$ws.Range("K70").Value = 5812.8
$ws.Range("L70").Value = 5862.5
$ws.Range("M70").Value = -5542.8
$ws.Range("N70").Value = -6402.5

$ws.Range("H73").Value = 5837.143
$ws.Range("I73").Value = 5812.8
$ws.Range("J73").Value = 5862.5
$ws.Range("K73").Value = 5812.8
$ws.Range("L73").Value = 5862.5
$ws.Range("M73").Value = -4876.8
$ws.Range("N73").Value = -7734.5

$ws.Range("H107").Value = 1690.2
$ws.Range("I107").Value = 1690.4
$ws.Range("J107").Value = 1690
$ws.Range("K107").Value = 1690.4
$ws.Range("L107").Value = 1690
$ws.Range("M107").Value = 229.5999999999999
$ws.Range("N107").Value = -5530

$ws.Range("H113").Value = 1655
$ws.Range("I113").Value = 1498
$ws.Range("J113").Value = 1742.2222
$ws.Range("K113").Value = 1498
$ws.Range("L113").Value = 1742.2222
$ws.Range("M113").Value = 672
$ws.Range("N113").Value = -6082.2222

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 450.41666
$ws.Range("I107").Value = 420.5
$ws.Range("K107").Value = 1261.5
$ws.Range("M107").Value = 658.5

$ws.Range("H136").Value = 3076.4443
$ws.Range("I136").Value = 2370.2222
$ws.Range("J136").Value = 4488.8887
$ws.Range("K136").Value = 7110.6666
$ws.Range("L136").Value = 13466.6661
$ws.Range("M136").Value = -4560.6666
$ws.Range("N136").Value = -18566.6661
